# Generate Report for Handback
# Updates localization-status.xlsx to reflect a completed handback:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#    (this is a shared string used by Overview!E/F and each language sheet's Status column,
#    so a single text update on one cell + reusing the same string elsewhere covers it)
#  - Each language sheet's row 2/3 gets its "Latest Target File" (I) and
#    "Latest Handback File" (J) populated, and "Latest Handback DateTime" (K) stamped
#  - New hyperlinks are added on the "Latest Target File" cells, mirroring the
#    existing hyperlink style/target used by column A
#  - Some columns are widened to fit the new (longer) text

$wb = $excel.ActiveWorkbook

$handbackStatus = "Handed back: in sync with en-US"

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cf8f43aa38dccff6ba8c4a7dba3034ff3d4c11c8/e2e/"

# ---- Overview sheet: widen the per-language status columns (E, F) ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $handbackStatus
$ov.Range("F2").Value = $handbackStatus
$ov.Range("E3").Value = $handbackStatus
$ov.Range("F3").Value = $handbackStatus
$ov.Columns.Item(5).ColumnWidth = 29.9777047293527
$ov.Columns.Item(6).ColumnWidth = 29.9777047293527

function Update-LangSheet {
    param(
        [string]$sheetName,
        [string]$file1,
        [string]$file2,
        [string]$handback1,
        [string]$handback2,
        [string]$dt1,
        [string]$dt2
    )

    $ws = $wb.Worksheets.Item($sheetName)

    # Status column (C) -> handed back
    $ws.Range("C2").Value = $handbackStatus
    $ws.Range("C3").Value = $handbackStatus

    # Latest Target File (I) - now populated with the source file name, hyperlinked
    $ws.Range("I2").Value = $file1
    $ws.Range("I2").Style = "HyperLink"
    $ws.Range("I3").Value = $file2
    $ws.Range("I3").Style = "HyperLink"

    # Latest Handback File (J)
    $ws.Range("J2").Value = $handback1
    $ws.Range("J3").Value = $handback2

    # Latest Handback DateTime (K)
    $ws.Range("K2").Value = $dt1
    $ws.Range("K3").Value = $dt2

    # Hyperlinks on the "Latest Target File" column, pointing at the same
    # source-file pages linked from column A.
    $ws.Hyperlinks.Add($ws.Range("I2"), ($ghBase + $file1), [Type]::Missing, [Type]::Missing, $file1)
    $ws.Hyperlinks.Add($ws.Range("I3"), ($ghBase + $file2), [Type]::Missing, [Type]::Missing, $file2)

    # Widen columns C, I, J to fit the new content
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}

# NOTE: this engine's PowerShell parser doesn't bind named (-param value)
# arguments to function params, so call positionally.
Update-LangSheet "zh-cn" `
    "06a1b3ac-3f75-4ab1-9ede-a483454965e9.md" `
    "e1366036-8bf9-40d2-a3d8-5229379de03f.md" `
    "06a1b3ac-3f75-4ab1-9ede-a483454965e9.a2b1bb8e90f482152ddfc670292d08168c31e606.zh-cn.xlf" `
    "e1366036-8bf9-40d2-a3d8-5229379de03f.3a408cc6bd88ad2c61d9c690b71865cb7d611df7.zh-cn.xlf" `
    "2016-08-17 04:55:28" `
    "2016-08-17 04:55:28"

Update-LangSheet "de-de" `
    "06a1b3ac-3f75-4ab1-9ede-a483454965e9.md" `
    "e1366036-8bf9-40d2-a3d8-5229379de03f.md" `
    "06a1b3ac-3f75-4ab1-9ede-a483454965e9.a2b1bb8e90f482152ddfc670292d08168c31e606.de-de.xlf" `
    "e1366036-8bf9-40d2-a3d8-5229379de03f.3a408cc6bd88ad2c61d9c690b71865cb7d611df7.de-de.xlf" `
    "2016-08-17 04:55:35" `
    "2016-08-17 04:55:35"

Write-Output "Handback report generated"
